# ------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet named "2022-Q1" right before the "总计"
#    (totals) sheet and populate it with the per-fund holdings for
#    that quarter (same layout as the other quarterly sheets).
# 2) Prepend a new row to the "总计" sheet summarising the 2022-Q1
#    totals (fund count + holding value), pushing the existing rows
#    down by one and keeping their running index in column A in sync.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Reference sheet whose header/row formatting we clone onto the new
# sheet so the visual style matches the other quarterly tabs exactly
# (copy/paste-format reuses the existing style ids instead of minting
# new ones).
$refSheet = $wb.Worksheets.Item("2021-Q4")

# ------------------------------------------------------------------
# 1) New "2022-Q1" sheet, inserted immediately before "总计"
# ------------------------------------------------------------------
$sheetTotal = $wb.Worksheets.Item("总计")
$ws = $wb.Worksheets.Add($sheetTotal)
$ws.Name = "2022-Q1"

# ---- header row -----------------------------------------------------
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# clone the bold / centered / bordered header look from the reference sheet
$refSheet.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

# ---- data rows --------------------------------------------------------
# columns: idx, code, name, size, stockPosition, positionRatio, value, rank
$fundRows = @(
    @(0, "166005", "中欧价值发现混合 -A",             "43.52", "93.97", "4.22", "1.8365", 4),
    @(1, "001882", "中欧价值发现混合 -E",             "43.52", "93.97", "4.22", "1.8365", 4),
    @(2, "001810", "中欧潜力价值灵活配置混合A",         "28.67", "94.05", "3.47", "0.9948", 6),
    @(3, "004232", "中欧价值发现混合 -C",             "10.98", "93.97", "4.22", "0.4634", 4),
    @(4, "166024", "中欧恒利三年定期开放混合",          "4.48",  "98.71", "5.39", "0.2415", 3),
    @(5, "005764", "中欧潜力价值灵活配置混合C",         "3.43",  "94.05", "3.47", "0.1190", 6),
    @(6, "001891", "中欧成长优选回报灵活配置混合E",      "2.97",  "94.42", "2.91", "0.0864", 7),
    @(7, "166020", "中欧成长优选回报灵活配置混合A",      "2.97",  "94.42", "2.91", "0.0864", 7),
    @(8, "009500", "国寿安保高股息混合A",              "0.96",  "73.47", "3.04", "0.0292", 10),
    @(9, "009501", "国寿安保高股息混合C",              "0.03",  "73.47", "3.04", "0.0009", 10)
)

$r = 2
foreach ($row in $fundRows) {
    $ws.Range("A$r").Value = $row[0]

    # fund code / size / position columns are kept as TEXT (leading
    # zeros in the fund codes must be preserved)
    $ws.Range("B$r").NumberFormat = "@"
    $ws.Range("B$r").Value = $row[1]

    $ws.Range("C$r").Value = $row[2]

    $ws.Range("D$r").NumberFormat = "@"
    $ws.Range("D$r").Value = $row[3]

    $ws.Range("E$r").NumberFormat = "@"
    $ws.Range("E$r").Value = $row[4]

    $ws.Range("F$r").NumberFormat = "@"
    $ws.Range("F$r").Value = $row[5]

    $ws.Range("G$r").NumberFormat = "@"
    $ws.Range("G$r").Value = $row[6]

    $ws.Range("H$r").Value = $row[7]

    $r = $r + 1
}

# clone the A-column (bold / centered / bordered index) formatting
$refSheet.Range("A2").Copy()
$ws.Range("A2:A11").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2) Add the 2022-Q1 summary row at the top of the "总计" sheet
#    (re-fetch the sheet handle - the one grabbed before Worksheets.Add
#    now points at the freshly inserted sheet instead)
# ------------------------------------------------------------------
$sheetTotal = $wb.Worksheets.Item("总计")
$sheetTotal.Rows.Item(2).Insert()
$sheetTotal.Range("B2:D2").ClearFormats()

$sheetTotal.Range("A2").Value = 0
$sheetTotal.Range("B2").Value = "2022-Q1"
$sheetTotal.Range("C2").Value = 10
$sheetTotal.Range("D2").Value = 5.69

# clone column-A formatting (style used by the other index cells, e.g. A3)
$sheetTotal.Range("A3").Copy()
$sheetTotal.Range("A2").PasteSpecial(-4122)
$sheetTotal.Range("A2").Value = 0

# re-sync the running index (column A) of the rows that got pushed
# down by the insert above
$sheetTotal.Range("A3").Value = 1
$sheetTotal.Range("A4").Value = 2
$sheetTotal.Range("A5").Value = 3
$sheetTotal.Range("A6").Value = 4
$sheetTotal.Range("A7").Value = 5

Write-Output "2022-Q1 sheet added and total sheet updated"
